# Auto-generated Excel COM-interop script to apply the Jenova_Profits.xlsx diff
# Updates cell values across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 812.6667
$ws.Range("J6").Value2 = 225
$ws.Range("L6").Value2 = 675
$ws.Range("N6").Value2 = -899
$ws.Range("H17").Value2 = 1752.75
$ws.Range("J17").Value2 = 1752.75
$ws.Range("L17").Value2 = 5258.25
$ws.Range("N17").Value2 = -5594.25
$ws.Range("H19").Value2 = 11112590
$ws.Range("I19").Value2 = 1144
$ws.Range("K19").Value2 = 1144
$ws.Range("M19").Value2 = -969
$ws.Range("H29").Value2 = 4159.8
$ws.Range("I29").Value2 = 650
$ws.Range("J29").Value2 = 6499.6665
$ws.Range("K29").Value2 = 1950
$ws.Range("L29").Value2 = 19498.9995
$ws.Range("M29").Value2 = -1669
$ws.Range("N29").Value2 = -20060.9995
$ws.Range("H132").Value2 = 6383.207
$ws.Range("I132").Value2 = 6383.207
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 19149.621
$ws.Range("L132").Value2 = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value2 = -16619.621
$ws.Range("H133").Value2 = 49339.184
$ws.Range("J133").Value2 = 49339.184
$ws.Range("L133").Value2 = 49339.184
$ws.Range("N133").Value2 = -59459.184
$ws.Range("H135").Value2 = 3416.375
$ws.Range("I135").Value2 = 3252.5862
$ws.Range("J135").Value2 = 4999.6665
$ws.Range("K135").Value2 = 29273.2758
$ws.Range("L135").Value2 = 44996.9985
$ws.Range("M135").Value2 = -26738.2758
$ws.Range("N135").Value2 = -50066.9985
$ws.Range("H138").Value2 = 6777.0723
$ws.Range("J138").Value2 = 7975.1665
$ws.Range("L138").Value2 = 23925.4995
$ws.Range("N138").Value2 = -34205.49950000001
$ws.Range("H141").Value2 = 2542
$ws.Range("I141").Value2 = 2491.111
$ws.Range("J141").Value2 = 3000
$ws.Range("K141").Value2 = 7473.333
$ws.Range("L141").Value2 = 9000
$ws.Range("M141").Value2 = -2293.333
$ws.Range("N141").Value2 = -19360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value2 = 9.800000000000001
$ws.Range("I3").Value2 = 9.800000000000001
$ws.Range("K3").Value2 = 9.800000000000001
$ws.Range("M3").Value2 = 105.2
$ws.Range("H61").Value2 = 3451.1
$ws.Range("I61").Value2 = 3451.1
$ws.Range("K61").Value2 = 3451.1
$ws.Range("M61").Value2 = -3239.1
$ws.Range("H74").Value2 = 310984.47
$ws.Range("I74").Value2 = 335368.97
$ws.Range("K74").Value2 = 335368.97
$ws.Range("M74").Value2 = -334494.97
$ws.Range("H77").Value2 = 310984.47
$ws.Range("I77").Value2 = 335368.97
$ws.Range("K77").Value2 = 1676844.85
$ws.Range("M77").Value2 = -1672476.85
$ws.Range("H102").Value2 = 2050.6667
$ws.Range("I102").Value2 = 2050.6667
$ws.Range("K102").Value2 = 2050.6667
$ws.Range("M102").Value2 = -428.6667000000002
$ws.Range("H132").Value2 = 186698.25
$ws.Range("I132").Value2 = 225879.11
$ws.Range("J132").Value2 = 15727.272
$ws.Range("K132").Value2 = 677637.33
$ws.Range("L132").Value2 = 47181.81600000001
$ws.Range("M132").Value2 = -675107.33
$ws.Range("N132").Value2 = -52241.81600000001
$ws.Range("H136").Value2 = 3451.1
$ws.Range("I136").Value2 = 3451.1
$ws.Range("K136").Value2 = 10353.3
$ws.Range("M136").Value2 = -7803.299999999999
$ws.Range("H137").Value2 = 0
$ws.Range("J137").Value2 = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value2 = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value2 = 1801.9
$ws.Range("I5").Value2 = 1254
$ws.Range("J5").Value2 = 2167.1667
$ws.Range("K5").Value2 = 1254
$ws.Range("L5").Value2 = 2167.1667
$ws.Range("M5").Value2 = -1141
$ws.Range("N5").Value2 = -2393.1667
$ws.Range("H11").Value2 = 397.5
$ws.Range("I11").Value2 = 183.33333
$ws.Range("J11").Value2 = 526
$ws.Range("K11").Value2 = 183.33333
$ws.Range("L11").Value2 = 526
$ws.Range("M11").Value2 = -43.33332999999999
$ws.Range("N11").Value2 = -806
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value2 = 0
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value2 = 0
$ws.Range("H86").Value2 = 1065226.8
$ws.Range("I86").Value2 = 1216291.1
$ws.Range("J86").Value2 = 7777
$ws.Range("K86").Value2 = 1216291.1
$ws.Range("L86").Value2 = 7777
$ws.Range("M86").Value2 = -1215168.1
$ws.Range("N86").Value2 = -10023
$ws.Range("H89").Value2 = 1065226.8
$ws.Range("I89").Value2 = 1216291.1
$ws.Range("J89").Value2 = 7777
$ws.Range("K89").Value2 = 6081455.5
$ws.Range("L89").Value2 = 38885
$ws.Range("M89").Value2 = -6075839.5
$ws.Range("N89").Value2 = -50117
$ws.Range("H105").Value2 = 6947638
$ws.Range("I105").Value2 = 1996.8889
$ws.Range("J105").Value2 = 11115023
$ws.Range("K105").Value2 = 1996.8889
$ws.Range("L105").Value2 = 11115023
$ws.Range("M105").Value2 = -249.8888999999999
$ws.Range("N105").Value2 = -11118517
$ws.Range("H132").Value2 = 50978.91
$ws.Range("J132").Value2 = 50978.91
$ws.Range("L132").Value2 = 50978.91
$ws.Range("N132").Value2 = -61098.91
$ws.Range("H135").Value2 = 63888.332
$ws.Range("J135").Value2 = 63888.332
$ws.Range("L135").Value2 = 63888.332
$ws.Range("N135").Value2 = -74028.33199999999
$ws.Range("H137").Value2 = 54408.59
$ws.Range("J137").Value2 = 54408.59
$ws.Range("L137").Value2 = 54408.59
$ws.Range("N137").Value2 = -64608.59
$ws.Range("H140").Value2 = 70963.336
$ws.Range("J140").Value2 = 70963.336
$ws.Range("L140").Value2 = 70963.336
$ws.Range("N140").Value2 = -81323.336

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 11433929
$ws.Range("J4").Value2 = 40010000
$ws.Range("L4").Value2 = 40010000
$ws.Range("N4").Value2 = -40010224
$ws.Range("H31").Value2 = 349396.1
$ws.Range("I31").Value2 = 738136
$ws.Range("J31").Value2 = 18247.297
$ws.Range("K31").Value2 = 738136
$ws.Range("L31").Value2 = 18247.297
$ws.Range("M31").Value2 = -737841
$ws.Range("N31").Value2 = -18837.297
$ws.Range("H34").Value2 = 349396.1
$ws.Range("I34").Value2 = 738136
$ws.Range("J34").Value2 = 18247.297
$ws.Range("K34").Value2 = 738136
$ws.Range("L34").Value2 = 18247.297
$ws.Range("M34").Value2 = -737934
$ws.Range("N34").Value2 = -18651.297
$ws.Range("H116").Value2 = 62124.75
$ws.Range("J116").Value2 = 62124.75
$ws.Range("L116").Value2 = 62124.75
$ws.Range("N116").Value2 = -71302.75
$ws.Range("H132").Value2 = 3311.3684
$ws.Range("J132").Value2 = 8331.666999999999
$ws.Range("L132").Value2 = 24995.001
$ws.Range("N132").Value2 = -30055.001
$ws.Range("H134").Value2 = 217900.58
$ws.Range("J134").Value2 = 781531.5600000001
$ws.Range("L134").Value2 = 2344594.68
$ws.Range("N134").Value2 = -2349664.68

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value2 = 685.875
$ws.Range("J92").Value2 = 719.5714
$ws.Range("L92").Value2 = 2158.7142
$ws.Range("N92").Value2 = -4654.7142
$ws.Range("H131").Value2 = 59685.688
$ws.Range("I131").Value2 = 46362.637
$ws.Range("J131").Value2 = 82232.38
$ws.Range("K131").Value2 = 139087.911
$ws.Range("L131").Value2 = 246697.14
$ws.Range("M131").Value2 = -134047.911
$ws.Range("N131").Value2 = -256777.14
$ws.Range("H132").Value2 = 2693238.5
$ws.Range("I132").Value2 = 9092188
$ws.Range("J132").Value2 = 27009.666
$ws.Range("K132").Value2 = 81829692
$ws.Range("L132").Value2 = 243086.994
$ws.Range("M132").Value2 = -81827162
$ws.Range("N132").Value2 = -248146.994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 359944.06
$ws.Range("I122").Value2 = 444121.84
$ws.Range("K122").Value2 = 1332365.52
$ws.Range("M122").Value2 = -1329915.52

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 5811.091
$ws.Range("I61").Value2 = 3601.2778
$ws.Range("J61").Value2 = 8462.866
$ws.Range("K61").Value2 = 3601.2778
$ws.Range("L61").Value2 = 8462.866
$ws.Range("M61").Value2 = -3399.2778
$ws.Range("N61").Value2 = -8866.866
$ws.Range("H82").Value2 = 2098.1177
$ws.Range("I82").Value2 = 1744.3636
$ws.Range("J82").Value2 = 2746.6667
$ws.Range("K82").Value2 = 1744.3636
$ws.Range("L82").Value2 = 2746.6667
$ws.Range("M82").Value2 = -1383.3636
$ws.Range("N82").Value2 = -3468.6667
$ws.Range("H85").Value2 = 2098.1177
$ws.Range("I85").Value2 = 1744.3636
$ws.Range("J85").Value2 = 2746.6667
$ws.Range("K85").Value2 = 1744.3636
$ws.Range("L85").Value2 = 2746.6667
$ws.Range("M85").Value2 = -496.3635999999999
$ws.Range("N85").Value2 = -5242.6667
$ws.Range("H113").Value2 = 5811.091
$ws.Range("I113").Value2 = 3601.2778
$ws.Range("J113").Value2 = 8462.866
$ws.Range("K113").Value2 = 3601.2778
$ws.Range("L113").Value2 = 8462.866
$ws.Range("M113").Value2 = -1431.2778
$ws.Range("N113").Value2 = -12802.866
$ws.Range("H132").Value2 = 6736.2144
$ws.Range("I132").Value2 = 5378
$ws.Range("K132").Value2 = 16134
$ws.Range("M132").Value2 = -13604

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 3911.4062
$ws.Range("J122").Value2 = 5818.091
$ws.Range("L122").Value2 = 17454.273
$ws.Range("N122").Value2 = -22354.273
$ws.Range("H126").Value2 = 5978.049
$ws.Range("I126").Value2 = 2463.1738
$ws.Range("K126").Value2 = 7389.5214
$ws.Range("M126").Value2 = -4919.5214
$ws.Range("H132").Value2 = 31253.25
$ws.Range("I132").Value2 = 3237.5667
$ws.Range("K132").Value2 = 9712.7001
$ws.Range("M132").Value2 = -7182.7001
